# LBRANDS_Q2_2013.xlsx — "chandra manual annot complete"
# Fills in the manual annotation scores (columns E:J — Clear, Assertive,
# Cautious, Optimistic, Specific, Relevant) for rows 2-34 of Sheet1, and
# updates the sheet view (zoom, frozen header row, final selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row (2-34) -> E,F,G,H,I,J scores, taken from the annotation diff.
$data = @{
    2  = @(2,2,1,2,1,2)
    3  = @(2,1,2,2,1,2)
    4  = @(2,2,1,2,2,2)
    5  = @(2,2,1,1,2,2)
    6  = @(2,2,1,2,2,2)
    7  = @(2,2,1,1,2,2)
    8  = @(2,2,1,1,2,2)
    9  = @(2,2,1,1,2,2)
    10 = @(2,2,2,1,1,2)
    11 = @(2,2,2,1,1,2)
    12 = @(2,2,2,2,2,2)
    13 = @(2,1,2,1,2,2)
    14 = @(2,2,1,1,1,2)
    15 = @(2,2,1,2,2,2)
    16 = @(2,1,1,1,1,2)
    17 = @(2,2,2,2,1,2)
    18 = @(2,2,1,2,2,2)
    19 = @(2,1,1,1,1,2)
    20 = @(2,2,1,2,1,2)
    21 = @(2,2,1,1,1,2)
    22 = @(2,2,1,1,1,2)
    23 = @(2,1,1,1,1,2)
    24 = @(2,2,1,2,1,2)
    25 = @(2,2,1,2,1,2)
    26 = @(2,2,2,2,1,2)
    27 = @(2,2,1,1,1,2)
    28 = @(2,1,1,1,1,2)
    29 = @(2,1,1,1,2,2)
    30 = @(2,2,1,1,1,2)
    31 = @(2,2,1,1,1,2)
    32 = @(2,1,1,1,1,2)
    33 = @(2,1,1,1,2,2)
    34 = @(2,2,1,1,2,2)
}

# Columns E..J are columns 5..10
$firstCol = 5

foreach ($rowNum in ($data.Keys | Sort-Object)) {
    $scores = $data[$rowNum]
    for ($i = 0; $i -lt $scores.Length; $i++) {
        $ws.Cells.Item($rowNum, $firstCol + $i).Value = $scores[$i]
    }
}

# --- Sheet view changes -------------------------------------------------
# Freeze the header row (row 1) in place.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Zoom to 85%.
$excel.ActiveWindow.Zoom = 85

# Final selection / scroll position, per the saved view state.
$ws.Range("E34").Select()
